$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 713.7
$ws.Range("C2").Value = 76

$ws.Range("B3").Value = 120
$ws.Range("C3").Value = 76

$ws.Range("B5").Value = 1188
$ws.Range("C5").Value = 152

$ws.Range("B6").Value = 131
$ws.Range("C6").Value = 76

$ws.Range("B8").Value = 508
$ws.Range("C8").Value = 76

$ws.Range("B9").Value = 77
$ws.Range("C9").Value = 38

$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 38

$ws.Range("B11").Value = 47
$ws.Range("C11").Value = 76
